$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G7").Value = "2016-08-17 18:40:54"
$wsZhCn.Range("H7").Value = "2016-08-17 18:40:49"
$wsDeDe.Range("H7").Value = "2016-08-17 18:40:54"
